$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3 and 4 used to store the dates as text ("2nd May" / "3rd May").
# Replace them with real date serial values so they match row 5's style
# (a proper date, not a shared-string label).
$ws.Range("B3").Value = 45048
$ws.Range("B4").Value = 45049

# Apply a custom date number format to the whole Date column so all the
# date cells (B1's header through B5) render the same way.
$ws.Columns.Item(2).NumberFormat = "mm/dd/yy;@"

# Give column B an explicit width, matching the new column formatting.
$ws.Columns.Item(2).ColumnWidth = 8.7265625

# Move the active selection.
$ws.Range("B7").Select()

# Switch the page to portrait orientation.
$ws.PageSetup.Orientation = 1

$wb.Save()
